$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 3518.75
$ws.Cells.Item(64, 9).Value = 3258.3333
$ws.Cells.Item(64, 10).Value = 3779.1667
$ws.Cells.Item(64, 11).Value = 3258.3333
$ws.Cells.Item(64, 12).Value = 3779.1667
$ws.Cells.Item(64, 13).Value = -3010.3333
$ws.Cells.Item(64, 14).Value = -4275.1667
$ws.Cells.Item(67, 8).Value = 3518.75
$ws.Cells.Item(67, 9).Value = 3258.3333
$ws.Cells.Item(67, 10).Value = 3779.1667
$ws.Cells.Item(67, 11).Value = 3258.3333
$ws.Cells.Item(67, 12).Value = 3779.1667
$ws.Cells.Item(67, 13).Value = -2400.3333
$ws.Cells.Item(67, 14).Value = -5495.1667
$ws.Cells.Item(96, 8).Value = 25001378
$ws.Cells.Item(96, 9).Value = 41667532
$ws.Cells.Item(96, 10).Value = 2144.25
$ws.Cells.Item(96, 11).Value = 125002596
$ws.Cells.Item(96, 12).Value = 6432.75
$ws.Cells.Item(96, 13).Value = -125001223
$ws.Cells.Item(96, 14).Value = -9178.75
$ws.Cells.Item(100, 8).Value = 2649.6428
$ws.Cells.Item(100, 10).Value = 3066.6667
$ws.Cells.Item(100, 12).Value = 3066.6667
$ws.Cells.Item(100, 14).Value = -4148.6667
$ws.Cells.Item(132, 8).Value = 3523.3333
$ws.Cells.Item(132, 9).Value = 3523.3333
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 10569.9999
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = -8039.999899999999
$ws.Cells.Item(132, 14).ClearContents()
$ws.Cells.Item(137, 8).Value = 296516.88
$ws.Cells.Item(137, 9).Value = 446921.78
$ws.Cells.Item(137, 11).Value = 1340765.34
$ws.Cells.Item(137, 13).Value = -1338215.34
$ws.Cells.Item(138, 8).Value = 2299.9783
$ws.Cells.Item(138, 10).Value = 3448.8572
$ws.Cells.Item(138, 12).Value = 10346.5716
$ws.Cells.Item(138, 14).Value = -20626.5716
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 24399.4
$ws.Cells.Item(32, 9).Value = 25597
$ws.Cells.Item(32, 11).Value = 25597
$ws.Cells.Item(32, 13).Value = -25310
$ws.Cells.Item(45, 8).Value = 2839.6052
$ws.Cells.Item(45, 9).Value = 2116.2778
$ws.Cells.Item(45, 10).Value = 3490.6
$ws.Cells.Item(45, 11).Value = 2116.2778
$ws.Cells.Item(45, 12).Value = 3490.6
$ws.Cells.Item(45, 13).Value = -1739.2778
$ws.Cells.Item(45, 14).Value = -4244.6
$ws.Cells.Item(61, 8).Value = 5099.8
$ws.Cells.Item(61, 9).Value = 3000
$ws.Cells.Item(61, 10).Value = 5333.1113
$ws.Cells.Item(61, 11).Value = 3000
$ws.Cells.Item(61, 12).Value = 5333.1113
$ws.Cells.Item(61, 13).Value = -2788
$ws.Cells.Item(61, 14).Value = -5757.1113
$ws.Cells.Item(74, 8).Value = 2108.05
$ws.Cells.Item(74, 9).Value = 1721.6875
$ws.Cells.Item(74, 11).Value = 1721.6875
$ws.Cells.Item(74, 13).Value = -847.6875
$ws.Cells.Item(77, 8).Value = 2108.05
$ws.Cells.Item(77, 9).Value = 1721.6875
$ws.Cells.Item(77, 11).Value = 8608.4375
$ws.Cells.Item(77, 13).Value = -4240.4375
$ws.Cells.Item(109, 8).Value = 37000
$ws.Cells.Item(109, 10).Value = 37000
$ws.Cells.Item(109, 12).Value = 37000
$ws.Cells.Item(109, 14).Value = -39774
$ws.Cells.Item(110, 8).Value = 4382.4443
$ws.Cells.Item(110, 9).Value = 4300.4
$ws.Cells.Item(110, 10).Value = 4485
$ws.Cells.Item(110, 11).Value = 4300.4
$ws.Cells.Item(110, 12).Value = 4485
$ws.Cells.Item(110, 13).Value = -2255.4
$ws.Cells.Item(110, 14).Value = -8575
$ws.Cells.Item(132, 8).Value = 25581.637
$ws.Cells.Item(132, 9).Value = 2151.818
$ws.Cells.Item(132, 10).Value = 49011.453
$ws.Cells.Item(132, 11).Value = 6455.454000000001
$ws.Cells.Item(132, 12).Value = 147034.359
$ws.Cells.Item(132, 13).Value = -3925.454000000001
$ws.Cells.Item(132, 14).Value = -152094.359
$ws.Cells.Item(136, 8).Value = 5099.8
$ws.Cells.Item(136, 9).Value = 3000
$ws.Cells.Item(136, 10).Value = 5333.1113
$ws.Cells.Item(136, 11).Value = 9000
$ws.Cells.Item(136, 12).Value = 15999.3339
$ws.Cells.Item(136, 13).Value = -6450
$ws.Cells.Item(136, 14).Value = -21099.3339
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 2244.923
$ws.Cells.Item(99, 9).Value = 1813.25
$ws.Cells.Item(99, 10).Value = 2935.6
$ws.Cells.Item(99, 11).Value = 1813.25
$ws.Cells.Item(99, 12).Value = 2935.6
$ws.Cells.Item(99, 13).Value = -315.25
$ws.Cells.Item(99, 14).Value = -5931.6
$ws.Cells.Item(105, 8).Value = 7145743
$ws.Cells.Item(105, 9).Value = 2950
$ws.Cells.Item(105, 10).Value = 16669467
$ws.Cells.Item(105, 11).Value = 2950
$ws.Cells.Item(105, 12).Value = 16669467
$ws.Cells.Item(105, 13).Value = -1203
$ws.Cells.Item(105, 14).Value = -16672961
$ws.Cells.Item(107, 8).Value = 1161.1666
$ws.Cells.Item(107, 9).Value = 769.2308
$ws.Cells.Item(107, 11).Value = 769.2308
$ws.Cells.Item(107, 13).Value = 1150.7692
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 144.75
$ws.Cells.Item(7, 9).Value = 150
$ws.Cells.Item(7, 10).Value = 139.5
$ws.Cells.Item(7, 11).Value = 150
$ws.Cells.Item(7, 12).Value = 139.5
$ws.Cells.Item(7, 13).Value = -37
$ws.Cells.Item(7, 14).Value = -365.5
$ws.Cells.Item(22, 8).Value = 462.5
$ws.Cells.Item(22, 10).Value = 200
$ws.Cells.Item(22, 12).Value = 200
$ws.Cells.Item(22, 14).Value = -900
$ws.Cells.Item(31, 8).Value = 14217.235
$ws.Cells.Item(31, 9).Value = 23493.941
$ws.Cells.Item(31, 11).Value = 23493.941
$ws.Cells.Item(31, 13).Value = -23198.941
$ws.Cells.Item(34, 8).Value = 14217.235
$ws.Cells.Item(34, 9).Value = 23493.941
$ws.Cells.Item(34, 11).Value = 23493.941
$ws.Cells.Item(34, 13).Value = -23291.941
$ws.Cells.Item(132, 8).Value = 20114.035
$ws.Cells.Item(132, 9).Value = 27446.158
$ws.Cells.Item(132, 10).Value = 4635.1113
$ws.Cells.Item(132, 11).Value = 82338.474
$ws.Cells.Item(132, 12).Value = 13905.3339
$ws.Cells.Item(132, 13).Value = -79808.474
$ws.Cells.Item(132, 14).Value = -18965.3339
$ws.Cells.Item(134, 8).Value = 1192
$ws.Cells.Item(134, 9).Value = 999.55554
$ws.Cells.Item(134, 11).Value = 2998.66662
$ws.Cells.Item(134, 13).Value = -463.66662
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(55, 8).Value = 2715.1
$ws.Cells.Item(55, 10).Value = 2715.1
$ws.Cells.Item(55, 12).Value = 8145.299999999999
$ws.Cells.Item(55, 14).Value = -8499.299999999999
$ws.Cells.Item(92, 8).Value = 1320
$ws.Cells.Item(92, 9).Value = 650
$ws.Cells.Item(92, 11).Value = 1950
$ws.Cells.Item(92, 13).Value = -702
$ws.Cells.Item(122, 8).Value = 813.4
$ws.Cells.Item(122, 10).Value = 966
$ws.Cells.Item(122, 12).Value = 8694
$ws.Cells.Item(122, 14).Value = -13594
$ws.Cells.Item(131, 8).Value = 760.5599999999999
$ws.Cells.Item(131, 10).Value = 765.5714
$ws.Cells.Item(131, 12).Value = 2296.7142
$ws.Cells.Item(131, 14).Value = -12376.7142
$ws.Cells.Item(138, 8).Value = 1802.5555
$ws.Cells.Item(138, 10).Value = 2077.6667
$ws.Cells.Item(138, 12).Value = 6233.000100000001
$ws.Cells.Item(138, 14).Value = -16513.0001
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 2660.7334
$ws.Cells.Item(113, 9).Value = 2090.2222
$ws.Cells.Item(113, 10).Value = 3516.5
$ws.Cells.Item(113, 11).Value = 2090.2222
$ws.Cells.Item(113, 12).Value = 3516.5
$ws.Cells.Item(113, 13).Value = 79.77779999999984
$ws.Cells.Item(113, 14).Value = -7856.5
$ws.Cells.Item(122, 8).Value = 1941.3636
$ws.Cells.Item(122, 9).Value = 1935.5
$ws.Cells.Item(122, 11).Value = 5806.5
$ws.Cells.Item(122, 13).Value = -3356.5
$ws.Cells.Item(132, 8).Value = 142110.19
$ws.Cells.Item(132, 9).Value = 149459
$ws.Cells.Item(132, 10).Value = 129249.75
$ws.Cells.Item(132, 11).Value = 448377
$ws.Cells.Item(132, 12).Value = 387749.25
$ws.Cells.Item(132, 13).Value = -445847
$ws.Cells.Item(132, 14).Value = -392809.25
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 4469.077
$ws.Cells.Item(7, 10).Value = 3633
$ws.Cells.Item(7, 12).Value = 3633
$ws.Cells.Item(7, 14).Value = -3857
$ws.Cells.Item(22, 8).Value = 3077.625
$ws.Cells.Item(22, 9).Value = 3815.3333
$ws.Cells.Item(22, 10).Value = 2635
$ws.Cells.Item(22, 11).Value = 3815.3333
$ws.Cells.Item(22, 12).Value = 2635
$ws.Cells.Item(22, 13).Value = -3520.3333
$ws.Cells.Item(22, 14).Value = -3225
$ws.Cells.Item(27, 8).Value = 3077.625
$ws.Cells.Item(27, 9).Value = 3815.3333
$ws.Cells.Item(27, 10).Value = 2635
$ws.Cells.Item(27, 11).Value = 3815.3333
$ws.Cells.Item(27, 12).Value = 2635
$ws.Cells.Item(27, 13).Value = -3708.3333
$ws.Cells.Item(27, 14).Value = -2849
$ws.Cells.Item(46, 8).Value = 1032.3636
$ws.Cells.Item(46, 9).Value = 829.375
$ws.Cells.Item(46, 10).Value = 1573.6666
$ws.Cells.Item(46, 11).Value = 829.375
$ws.Cells.Item(46, 12).Value = 1573.6666
$ws.Cells.Item(46, 13).Value = -641.375
$ws.Cells.Item(46, 14).Value = -1949.6666
$ws.Cells.Item(61, 8).Value = 4424.625
$ws.Cells.Item(61, 9).Value = 2199.125
$ws.Cells.Item(61, 10).Value = 8875.625
$ws.Cells.Item(61, 11).Value = 2199.125
$ws.Cells.Item(61, 12).Value = 8875.625
$ws.Cells.Item(61, 13).Value = -1997.125
$ws.Cells.Item(61, 14).Value = -9279.625
$ws.Cells.Item(109, 8).Value = 31196.25
$ws.Cells.Item(109, 10).Value = 31196.25
$ws.Cells.Item(109, 12).Value = 31196.25
$ws.Cells.Item(109, 14).Value = -33970.25
$ws.Cells.Item(113, 8).Value = 4424.625
$ws.Cells.Item(113, 9).Value = 2199.125
$ws.Cells.Item(113, 10).Value = 8875.625
$ws.Cells.Item(113, 11).Value = 2199.125
$ws.Cells.Item(113, 12).Value = 8875.625
$ws.Cells.Item(113, 13).Value = -29.125
$ws.Cells.Item(113, 14).Value = -13215.625
$ws.Cells.Item(123, 8).Value = 0
$ws.Cells.Item(123, 10).Value = 0
$ws.Cells.Item(123, 12).Value = 0
$ws.Cells.Item(123, 14).ClearContents()
$ws.Cells.Item(126, 8).Value = 4469.077
$ws.Cells.Item(126, 10).Value = 3633
$ws.Cells.Item(126, 12).Value = 10899
$ws.Cells.Item(126, 14).Value = -15839
$ws.Cells.Item(132, 8).Value = 367140.44
$ws.Cells.Item(132, 9).Value = 483785.47
$ws.Cells.Item(132, 10).Value = 2624.625
$ws.Cells.Item(132, 11).Value = 1451356.41
$ws.Cells.Item(132, 12).Value = 7873.875
$ws.Cells.Item(132, 13).Value = -1448826.41
$ws.Cells.Item(132, 14).Value = -12933.875
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(27, 8).Value = 29677.8
$ws.Cells.Item(27, 9).Value = 15000
$ws.Cells.Item(27, 10).Value = 31308.666
$ws.Cells.Item(27, 11).Value = 15000
$ws.Cells.Item(27, 12).Value = 31308.666
$ws.Cells.Item(27, 14).Value = -31446.666
$ws.Cells.Item(27, 13).Value = -14931
$ws.Cells.Item(39, 8).Value = 4980
$ws.Cells.Item(39, 9).Value = 0
$ws.Cells.Item(39, 10).Value = 4980
$ws.Cells.Item(39, 11).Value = 0
$ws.Cells.Item(39, 12).Value = 4980
$ws.Cells.Item(39, 13).ClearContents()
$ws.Cells.Item(39, 14).Value = -5806
$ws.Cells.Item(92, 8).Value = 18350
$ws.Cells.Item(92, 10).Value = 18350
$ws.Cells.Item(92, 12).Value = 18350
$ws.Cells.Item(92, 14).Value = -23342
$ws.Cells.Item(113, 8).Value = 1802718.5
$ws.Cells.Item(113, 9).Value = 1019.2308
$ws.Cells.Item(113, 11).Value = 3057.6924
$ws.Cells.Item(113, 13).Value = -887.6923999999999
$ws.Cells.Item(132, 8).Value = 2571.3914
$ws.Cells.Item(132, 9).Value = 2307
$ws.Cells.Item(132, 11).Value = 6921
$ws.Cells.Item(132, 13).Value = -4391
